$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 640
$ws.Range("I18").Value = 675
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 675
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -391
$ws.Range("N18").Value = -1068

$ws.Range("H113").Value = 71431570
$ws.Range("I113").Value = 125001250
$ws.Range("J113").Value = 5329.3335
$ws.Range("K113").Value = 125001250
$ws.Range("L113").Value = 5329.3335
$ws.Range("M113").Value = -124997996
$ws.Range("N113").Value = -11837.3335

$ws.Range("H129").Value = 820.3103599999999
$ws.Range("I129").Value = 405.73685
$ws.Range("J129").Value = 1608
$ws.Range("K129").Value = 1217.21055
$ws.Range("L129").Value = 4824
$ws.Range("M129").Value = 3782.78945
$ws.Range("N129").Value = -14824

$ws.Range("H137").Value = 3410277.5
$ws.Range("I137").Value = 1563678.5
$ws.Range("J137").Value = 8334541.5
$ws.Range("K137").Value = 4691035.5
$ws.Range("L137").Value = 25003624.5
$ws.Range("M137").Value = -4688485.5
$ws.Range("N137").Value = -25008724.5

$ws.Range("H138").Value = 1761.659
$ws.Range("I138").Value = 1155.2222
$ws.Range("J138").Value = 2181.5
$ws.Range("K138").Value = 3465.6666
$ws.Range("L138").Value = 6544.5
$ws.Range("M138").Value = 1674.3334
$ws.Range("N138").Value = -16824.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2723904.2
$ws.Range("I32").Value = 3845.825
$ws.Range("J32").Value = 20857628
$ws.Range("K32").Value = 3845.825
$ws.Range("L32").Value = 20857628
$ws.Range("M32").Value = -3558.825
$ws.Range("N32").Value = -20858202

$ws.Range("H45").Value = 2813.6086
$ws.Range("I45").Value = 1054.3077
$ws.Range("J45").Value = 5100.7
$ws.Range("K45").Value = 1054.3077
$ws.Range("L45").Value = 5100.7
$ws.Range("M45").Value = -677.3077000000001
$ws.Range("N45").Value = -5854.7

$ws.Range("H74").Value = 4762860.5
$ws.Range("I74").Value = 940.8158
$ws.Range("J74").Value = 50001096
$ws.Range("K74").Value = 940.8158
$ws.Range("L74").Value = 50001096
$ws.Range("M74").Value = -66.81579999999997
$ws.Range("N74").Value = -50002844

$ws.Range("H77").Value = 4762860.5
$ws.Range("I77").Value = 940.8158
$ws.Range("J77").Value = 50001096
$ws.Range("K77").Value = 4704.079
$ws.Range("L77").Value = 250005480
$ws.Range("M77").Value = -336.0789999999997
$ws.Range("N77").Value = -250014216

$ws.Range("H132").Value = 82524.92
$ws.Range("I132").Value = 99282.336
$ws.Range("J132").Value = 4831.4546
$ws.Range("K132").Value = 297847.008
$ws.Range("L132").Value = 14494.3638
$ws.Range("M132").Value = -295317.008
$ws.Range("N132").Value = -19554.3638

$ws.Range("H134").Value = 35120
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 35120
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 35120
$ws.Range("N134").Value = -45260

$ws.Range("H139").Value = 41571.5
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 43885.8
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 43885.8
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -54165.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 39962.676
$ws.Range("I134").Value = 47787.188
$ws.Range("J134").Value = 1442
$ws.Range("K134").Value = 143361.564
$ws.Range("L134").Value = 4326
$ws.Range("M134").Value = -140826.564
$ws.Range("N134").Value = -9396

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2700
$ws.Range("I31").Value = 1984.3334
$ws.Range("J31").Value = 3654.2222
$ws.Range("K31").Value = 1984.3334
$ws.Range("L31").Value = 3654.2222
$ws.Range("M31").Value = -1689.3334
$ws.Range("N31").Value = -4244.2222

$ws.Range("H34").Value = 2700
$ws.Range("I34").Value = 1984.3334
$ws.Range("J34").Value = 3654.2222
$ws.Range("K34").Value = 1984.3334
$ws.Range("L34").Value = 3654.2222
$ws.Range("M34").Value = -1782.3334
$ws.Range("N34").Value = -4058.2222

$ws.Range("H52").Value = 25093
$ws.Range("I52").Value = 23800
$ws.Range("J52").Value = 25277.715
$ws.Range("K52").Value = 23800
$ws.Range("L52").Value = 25277.715
$ws.Range("M52").Value = -23506
$ws.Range("N52").Value = -25865.715

$ws.Range("H132").Value = 2110.0454
$ws.Range("I132").Value = 1958.6
$ws.Range("J132").Value = 3624.5
$ws.Range("K132").Value = 5875.799999999999
$ws.Range("L132").Value = 10873.5
$ws.Range("M132").Value = -3345.799999999999
$ws.Range("N132").Value = -15933.5

$ws.Range("H134").Value = 4274.3447
$ws.Range("I134").Value = 4348.5557
$ws.Range("J134").Value = 3272.5
$ws.Range("K134").Value = 13045.6671
$ws.Range("L134").Value = 9817.5
$ws.Range("M134").Value = -10510.6671
$ws.Range("N134").Value = -14887.5

$ws.Range("H141").Value = 27758.666
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 27758.666
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 27758.666
$ws.Range("N141").Value = -38118.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1333938.8
$ws.Range("I5").Value = 592.3333
$ws.Range("J5").Value = 3333958.2
$ws.Range("K5").Value = 1776.9999
$ws.Range("L5").Value = 10001874.6
$ws.Range("M5").Value = -1664.9999
$ws.Range("N5").Value = -10002098.6

$ws.Range("H129").Value = 1575.0769
$ws.Range("I129").Value = 716.3570999999999
$ws.Range("J129").Value = 2576.9167
$ws.Range("K129").Value = 2149.0713
$ws.Range("L129").Value = 7730.750100000001
$ws.Range("M129").Value = 2850.9287
$ws.Range("N129").Value = -17730.7501

$ws.Range("H135").Value = 1333938.8
$ws.Range("I135").Value = 592.3333
$ws.Range("J135").Value = 3333958.2
$ws.Range("K135").Value = 5330.9997
$ws.Range("L135").Value = 30005623.8
$ws.Range("M135").Value = -2795.9997
$ws.Range("N135").Value = -30010693.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2948.0557
$ws.Range("I122").Value = 3175.111
$ws.Range("J122").Value = 2721
$ws.Range("K122").Value = 9525.332999999999
$ws.Range("L122").Value = 8163
$ws.Range("M122").Value = -7075.332999999999
$ws.Range("N122").Value = -13063

$ws.Range("H132").Value = 2185.818
$ws.Range("I132").Value = 1757.4286
$ws.Range("J132").Value = 2935.5
$ws.Range("K132").Value = 5272.2858
$ws.Range("L132").Value = 8806.5
$ws.Range("M132").Value = -2742.2858
$ws.Range("N132").Value = -13866.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2668.261
$ws.Range("I122").Value = 2152
$ws.Range("J122").Value = 3282.8572
$ws.Range("K122").Value = 6456
$ws.Range("L122").Value = 9848.571599999999
$ws.Range("M122").Value = -4006
$ws.Range("N122").Value = -14748.5716

$ws.Range("H132").Value = 2147.3784
$ws.Range("I132").Value = 1795.4688
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 5386.4064
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -2856.4064
$ws.Range("N132").Value = -18258.8

$ws.Range("H133").Value = 42848.43
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 42848.43
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 42848.43
$ws.Range("N133").Value = -47908.43

$ws.Range("H136").Value = 1350.329
$ws.Range("I136").Value = 1248.0164
$ws.Range("J136").Value = 1766.4
$ws.Range("K136").Value = 3744.0492
$ws.Range("L136").Value = 5299.200000000001
$ws.Range("M136").Value = -1194.0492
$ws.Range("N136").Value = -10399.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49888.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 49888.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 49888.75
$ws.Range("N46").Value = -50350.75

$ws.Range("H113").Value = 27054.87
$ws.Range("I113").Value = 40325.2
$ws.Range("J113").Value = 1535
$ws.Range("K113").Value = 120975.6
$ws.Range("L113").Value = 4605
$ws.Range("M113").Value = -118805.6
$ws.Range("N113").Value = -8945

$ws.Range("H122").Value = 3695.9412
$ws.Range("I122").Value = 2444.5715
$ws.Range("J122").Value = 4571.9
$ws.Range("K122").Value = 7333.7145
$ws.Range("L122").Value = 13715.7
$ws.Range("M122").Value = -4883.7145
$ws.Range("N122").Value = -18615.7

$ws.Range("H132").Value = 1246.9412
$ws.Range("I132").Value = 1121.3158
$ws.Range("J132").Value = 1897.909
$ws.Range("K132").Value = 3363.9474
$ws.Range("L132").Value = 5693.727000000001
$ws.Range("M132").Value = -833.9474
$ws.Range("N132").Value = -10753.727

$ws.Range("H134").Value = 49888.75
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 49888.75
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 149666.25
$ws.Range("N134").Value = -154736.25

$ws.Range("H136").Value = 1344.2063
$ws.Range("I136").Value = 1352.5927
$ws.Range("J136").Value = 1293.8889
$ws.Range("K136").Value = 4057.7781
$ws.Range("L136").Value = 3881.6667
$ws.Range("M136").Value = -1507.7781
$ws.Range("N136").Value = -8981.6667

Write-Output "edits applied"
